$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.912493824958801
$ws.Range("B1").Value = 3.451638698577881
$ws.Range("C1").Value = 2.903152942657471
$ws.Range("D1").Value = 2.053718566894531
$ws.Range("E1").Value = 1.191772937774658
